$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unlock only the specific cells we
# need to touch, edit them, then re-lock so the sheet protection state
# (password hash / flags) is left exactly as it was.

$targets = @("A7", "D2", "E2", "D3", "E3", "E4")
foreach ($addr in $targets) {
    $ws.Range($addr).Locked = $false
}

# Update the confidential disclaimer date from 2021-05-27 to 2021-05-28
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."
# Setting a taller wrapped value can auto-grow the row height; restore the
# original (default) row height so only the cell content changes.
$ws.Rows.Item(7).AutoFit()

# Update weight/percent-change values for rows 2-4
$ws.Range("D2").Value = 0.8472486831783033
$ws.Range("E2").Value = 0.002356150793650924

$ws.Range("D3").Value = 0.1527513168216968
$ws.Range("E3").Value = 0.006419662509170809

$ws.Range("E4").Value = 0.002976857559117141

foreach ($addr in $targets) {
    $ws.Range($addr).Locked = $true
}
